$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value updates derived from the authoritative diff
$updates = @{
    "G9" = 1.48
    "I9" = 6.25
    "W9" = 10
    "AF9" = 34
    "AG9" = 21
    "AI9" = 51
    "G20" = 2.18
    "H20" = 3.15
    "L20" = 1.45
    "M20" = 2.4
    "N20" = 2.27
    "O20" = 1.5
    "P20" = 1.47
    "Q20" = 2.32
    "R20" = 2.02
    "S20" = 1.62
    "T20" = 5.9
    "U20" = 9
    "V20" = 9.5
    "W20" = 20
    "X20" = 21
    "Y20" = 40
    "Z20" = 7.1
    "AA20" = 6.2
    "AE20" = 7.6
    "AG20" = 12
    "AH20" = 40
    "N28" = 2.02
    "T28" = 7.9
    "X28" = 23
    "Y28" = 32
    "AE28" = 8.75
    "AI28" = 23
    "AE29" = 9.25
    "AF29" = 21
    "J30" = 1.07
    "K30" = 6.7
    "L30" = 1.34
    "M30" = 3
    "N30" = 2
    "O30" = 1.72
    "Q30" = 2.52
    "R30" = 1.83
    "S30" = 1.87
    "T30" = 6.6
    "U30" = 8.75
    "V30" = 8.25
    "X30" = 15.5
    "Y30" = 28
    "Z30" = 6.7
    "AB30" = 15
    "AC30" = 75
    "AD30" = 600
    "AE30" = 10.75
    "AJ30" = 45
    "G32" = 1.62
    "H32" = 3.6
    "I32" = 5.5
    "R32" = 1.87
    "S32" = 1.85
    "T32" = 6.6
    "U32" = 7.5
    "V32" = 7.9
    "W32" = 12
    "X32" = 12.5
    "AA32" = 7
    "AB32" = 16.5
    "AC32" = 80
    "AD32" = 700
    "AE32" = 13
    "AF32" = 32
    "AG32" = 17.5
    "AH32" = 110
    "AI32" = 60
    "AJ32" = 60
    "T33" = 6.2
    "V33" = 9.75
    "X33" = 24
    "Y33" = 40
    "AA33" = 5.7
    "AC33" = 110
    "AE33" = 7.2
    "AF33" = 14.5
    "AG33" = 11.75
    "AI33" = 35
    "AJ33" = 50
    "G41" = 2.82
    "H41" = 3.35
    "I41" = 2.27
    "L41" = 1.26
    "M41" = 3.15
    "N41" = 1.78
    "O41" = 1.82
    "T41" = 9.75
    "U41" = 15
    "V41" = 10.25
    "W41" = 35
    "X41" = 23
    "Y41" = 30
    "Z41" = 10.75
    "AA41" = 6.6
    "AB41" = 13.5
    "AE41" = 8.5
    "AF41" = 11.5
    "AG41" = 9
    "AH41" = 23
    "AI41" = 18
    "AJ41" = 27
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
